# Edit Review_314.docx per commit diff: swap CLIP/localization review content
# for "Were RNNs All We Needed?" review content (new title, body paragraphs, link);
# drop the short closing remark paragraph and update the arxiv link.

$d = $word.ActiveDocument

# --- Paragraph 1: two runs separated by a manual line break (<w:br/>) ---
# Run 1: update the date in the headline
$d.Paragraphs.Item(1).Range.Find.Execute('-06.10.24:', $false, $false, $false, $false, $false, $true, 1, $false, '-04.10.24:', 2) | Out-Null
# Run 2: update the paper title
$d.Paragraphs.Item(1).Range.Find.Execute('CONTRASTIVE LOCALIZED LANGUAGE-IMAGE PRE-TRAINING', $false, $false, $false, $false, $false, $true, 1, $false, 'Were RNNs All We Needed?', 2) | Out-Null

# --- Paragraph 2 ---
# The new text contains a straight double-quote character. Word's Find/Replace
# "smart quotes" AutoFormat mangles a literal " inside ReplaceWith text, so we swap
# in a private-use placeholder character first and then patch just that one spot
# back to a literal " via a plain Range.Text assignment (bypasses AutoFormat).
$placeholder = [char]0xE000
$newP2 = 'המאמר הזה משך את תשומת ליבי כי יש לו ״all we needed' + $placeholder + ' בכותרת. מסיבה שאינה ב-100% ברורה לי מאמרים כאלו יוצרים בי דחף חזק לסקור אותם. אז ככה הגעתי למאמר הזה שאלולא השם כנראה שלא הייתי מגיע אליו.'
$d.Paragraphs.Item(2).Range.Find.Execute('ממשיכים הפסקה בסקירות על מודלי שפה ועוברים לסקירות על מודלים מולטימודליים (שפה ותמונות). טוב, הפסקה למחצה. אתם בטח זוכרים את המודל שנקרא CLIP שעשה הרבה רעש לפני כמה שנים. ', $false, $false, $false, $false, $false, $true, 1, $false, $newP2, 2) | Out-Null
$phFind = $d.Content
$phFind.Find.Execute($placeholder) | Out-Null
$d.Range($phFind.Start, $phFind.End).Text = '"'

# --- Paragraph 3 ---
$d.Paragraphs.Item(3).Range.Find.Execute('CLIP הוא אחד המודלים מולטימודליים הראשוניים שהצליח לייצר אמבדינגס חזקים ומיושרים (aligned) של טקסט ושל תמונות. מיושרים הכוונה של הייצוגים של תמונה וטקסט שמתאר את תוכנה קרובים אחד לשני בזמן שהייצוגים של תמונה וטקסט לא מתאימים רחוקים אחד מהשני (במקרה הזה ביחס למרחק קוסיין ביניהם).', $false, $false, $false, $false, $false, $true, 1, $false, 'המאמר מציע לשפצר את ה-RNN כך שנוכל להפעיל אותו בצורה מקבילית במהלך האימון. הסיבה העיקרית ש-RNN כמעט יצא מכלל שימוש היום הוא חוסר היכולת שלו להתאמן באופן מקבילי כלומר לבצע חיזוי של כמה טוקנים ממוסכים. הטרנספורמרים לעומת זאת כן ניחנים ביכולת הזו אך יש להם מגבלה בדמות סיבוכיות ריבועית במונחי אורך הסדרה (שכואבת לנו בעיקר באינפרנס כי מאמנים אותם פעם אחת) שמקשה על השימוש (לפחות הנאיבי שלהם) לסדרות מאוד ארוכות.', 2) | Out-Null

# --- Paragraph 4 ---
$d.Paragraphs.Item(4).Range.Find.Execute('המודל הזה אומן על דאטהסט ענק של תמונות והכותרות שלהם (או טאגים) מהאינטרנט כאשר אימנו אותו תוך שימוש בטכניקה למידה ניגודית (contrastive learning או CL). בגדול מאוד טכניקות CL מאומנות להפיק ייצוג סמנטי מדאטה (מסוגים שונים) כאשר המטרה היא לקרב את הייצוגים (אמבדינגס) של פיסות דאטה קרובות (או חיוביות) ולהרחיק ייצוגים של פיסות דאטה לא דומות (שליליות). במקרה של CLIP פיסות דאטה חיוביות הם הייצוגים של תמונה והכותרת שלה ואילו הזוגות השליליים בנויים מכותבות ותמונות שנבחרו באקראי.', $false, $false, $false, $false, $false, $true, 1, $false, 'מצד שני ל-RNNs יש יכולת יותר טובה לעבד סדרות מאוד ארוכות כי כל ה״זיכרון״ שלהם מקודד בכמה ווקטורים (1,2 או 3) והסיבוכיות החישובית שלהם פרופורציונלית לאורך הסדרה ולא לריבוע שלה (גם באימון וגם באינפרנס). כאמור הבעיה הגדולה של ה-RNNS שדי הרגה את הארכיטקטורה הזו היא אי יכולתה לאפשר חיזוי מקבילי באימון. זה שהופך את האימון על כמויות דאטה עצומות כמו שמקובל היום (עשרות טריליונים טוקנים) עם RNNs לארוך מדי ולא פיזיבילי. ', 2) | Out-Null

# --- Paragraph 5 ---
$d.Paragraphs.Item(5).Range.Find.Execute('המאמר שנסקור אחד כאמור משכלל את CLIP על ידי הקניה של יכולות לוקליזציה לייצוג. הכוונה כאן שהמחברים מאמנים ייצוגים של תמונה ושל טקסט באופן כזה שבהינתן ייצוג התמונה I וייצוג התיאור של פאץ'' ב I המכיל אובייקט מסוים יהיה ניתן להפיק ב״קלות״ את מיקום האובייקט בתמונה. ', $false, $false, $false, $false, $false, $true, 1, $false, 'חשוב להבין שהסיבה לחוסר יכולת לחזות בצורה מקבילי נובעת מהמעברים הלא לינאריים בין המצבים החבויים ב-RNN (גם ב-LSTM וגם ב-GRU). ', 2) | Out-Null

# --- Paragraph 6 ---
$d.Paragraphs.Item(6).Range.Find.Execute('במילים פשוטות נניח שיש לנו אריה עומד ושואג בתמונה הנמצא ב-bounding box (המוגדר על ידי רביעיה של קואורדינטות שלו בתמונה) המסומן ב- B. המחברים מאמנים רשת אנקודר לתמונות f_I רשת אנקודר לטקסט f_T כך שייצוג התמונה R_I ייצוג ״אריה עומד ושואג״ R_T, המופקים על ידי שני האנקודר האלו (בהתאמה) כך שרשת רדודה יחסית (נקראת prompter במאמר), המקבלת אותם, תוכל לחזות את מיקום האריה B בתמונה. דרך אגב המיקום כאן לא חייב להיות מתואר על ידי bounding box אלא יכול להיות מוגדר (בערך) על ידי כמה ניקודת, תיאור כללי (נגיד חיה, בלי להזכיר שזה אריה) ובעוד צורות.', $false, $false, $false, $false, $false, $true, 1, $false, 'לאחרונה SSMs (או State Space Models) ניסו לטפל בבעיה הזו דרך ארכיטקטורה שבה המעברים האלו כן לינאריים וארכיטקטורת ממבה (שסקרתי בהרחבה לפני כמה חודשים) ששכללה SSMs לרמת ביצועים קרובה לטרנספורמרים. בנוסף A21 Labs השתמשו בממבה כאבן בניין של הארכיטקטורה החדשה שלהם לפני כחודשיים(יחד עם הטרנספורמרים).', 2) | Out-Null

# --- Paragraph 7 ---
$d.Paragraphs.Item(7).Range.Find.Execute('האימון נעשה כמו בלמידה הניגודית כמו ב-CLIP המקורי. אבל בנוסף ללוס הרגיל שלו יש כאן עוד לוס ניגודי המקרב את ייצוגים של כותרת הפאץ'' בתמונה לייצוג המופק על Prompter מייצוג התמונה ומהמתאר של הפאץ'' (נגיד BB) ומרחיק את הייצוגים האלו לפאצ''ים שונים. כמובן שה-Prompter גם מאומן תוך כדי,', $false, $false, $false, $false, $false, $true, 1, $false, 'עכשיו אתם שואלים מה המאמר המסוקר עשה בנידון. כאמור הבעיה הגדולה ב-RNN היה מעברים לא לינאריים בין המצבים החבויים. המחברים פשוט הורידו את התלות הלא לינארית מהמשוואות של LSTM ו-GRU. מה שהתקבל כתוצאה מכך ניתן למקבול במהלך האימון (אבל דורש יותר זיכרון מהגרסאות הרגילות). יצא משהו די דומה לממבה - גם כן המצב החבוי תלוי באופן ליניארי במצב החבוי הקודם ובאופן לא לינארי בייצוג האיבר הנוכחי של סדרת הדאטה.', 2) | Out-Null

# --- Paragraph 8 ---
$d.Paragraphs.Item(8).Range.Find.Execute('המאמר משתמש במודלים מאומנים למטרת זיהוי אובייקטים בתמונה (OWLv2) ובמודלים מאומנים אחרים (VeCap) למתן כותרות לפאצ''ים האלו. ', $false, $false, $false, $false, $false, $true, 1, $false, 'מה שמפתיע אותי קצת כאן זה ביצועים טובים מדי - אני קצת חשדן אבל בואו נראה מה קורה עם הארכיטקטורה הזו בעתיד.', 2) | Out-Null

# --- Delete the short "מאמר די חמוד וקליל..." paragraph (old paragraph 9) entirely ---
# (this also pulls the final link paragraph up to become the new paragraph 9)
$d.Paragraphs.Item(9).Range.Delete()

# --- Paragraph 9 (now the arxiv-link paragraph) - update the URL ---
$d.Paragraphs.Item(9).Range.Find.Execute('https://arxiv.org/pdf/2410.02746', $false, $false, $false, $false, $false, $true, 1, $false, 'https://arxiv.org/abs/2410.01201v1', 2) | Out-Null

Write-Output "Edit complete"
